$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap Bolivia/Tunez order: row 90 becomes Bolivia, row 91 becomes Tunez
$ws.Range("A90").Value = "Bolivia"
$ws.Range("A91").Value = "Tunez"

# Update Bolivia row (90) with new stats
$ws.Cells.Item(90, 2).Value = 950
$ws.Cells.Item(90, 3).Value = 84
$ws.Cells.Item(90, 4).Value = 80
$ws.Cells.Item(90, 5).Value = 820
$ws.Cells.Item(90, 6).Value = 3
$ws.Cells.Item(90, 7).Value = 4
$ws.Cells.Item(90, 8).Value = 50

# Update Tunez row (91) with its (previous) stats
$ws.Cells.Item(91, 2).Value = 949
$ws.Cells.Item(91, 3).Value = 0
$ws.Cells.Item(91, 4).Value = 216
$ws.Cells.Item(91, 5).Value = 695
$ws.Cells.Item(91, 6).Value = 20
$ws.Cells.Item(91, 7).Value = 0
$ws.Cells.Item(91, 8).Value = 38

# China (row 12)
$ws.Cells.Item(12, 2).Value = 82830
$ws.Cells.Item(12, 3).Value = 3
$ws.Cells.Item(12, 4).Value = 77474
$ws.Cells.Item(12, 5).Value = 723
$ws.Cells.Item(12, 6).Value = 52
$ws.Cells.Item(12, 7).Value = 1
$ws.Cells.Item(12, 8).Value = 4633

# Mexico (row 28)
$ws.Cells.Item(28, 2).Value = 14677
$ws.Cells.Item(28, 3).Value = 835
$ws.Cells.Item(28, 4).Value = 8354
$ws.Cells.Item(28, 5).Value = 4972
$ws.Cells.Item(28, 6).Value = 378
$ws.Cells.Item(28, 7).Value = 46
$ws.Cells.Item(28, 8).Value = 1351

# Panama (row 49)
$ws.Cells.Item(49, 2).Value = 5779
$ws.Cells.Item(49, 3).Value = 241
$ws.Cells.Item(49, 4).Value = 338
$ws.Cells.Item(49, 5).Value = 5276
$ws.Cells.Item(49, 6).Value = 85
$ws.Cells.Item(49, 7).Value = 6
$ws.Cells.Item(49, 8).Value = 165

# Uruguay (row 104)
$ws.Cells.Item(104, 2).Value = 606
$ws.Cells.Item(104, 3).Value = 10
$ws.Cells.Item(104, 4).Value = 375
$ws.Cells.Item(104, 5).Value = 216
$ws.Cells.Item(104, 6).Value = 10
$ws.Cells.Item(104, 7).Value = 1
$ws.Cells.Item(104, 8).Value = 15
